# "user schema id 변경" — the users table schema (rows 2-3, columns D-L) is
# updated: two new columns (memberEmail / memberType) are appended, and the
# descriptions for the existing "id" / "name" columns are updated to reflect
# that they now come from the SSO login token. The old temporary note about
# memberID/memberPW (previously in L2) is removed since it is now obsolete.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New columns appended to the "users" table schema
$ws.Range("J2").Value = "memberEmail"
$ws.Range("K2").Value = "memberType"
$ws.Range("J3").Value = "token에 있는 정보"
$ws.Range("K3").Value = "token에 있는 정보"

# "id" / "name" columns now sourced from the SSO token
$ws.Range("F3").Value = "식별자(token에 있는 정보)"
$ws.Range("G3").Value = "token에 있는 정보"

# Obsolete temporary note about memberID/memberPW is removed
$ws.Range("L2").Value = ""

# Move the sheet selection to reflect where the edit was made
[void]$ws.Activate()
[void]$ws.Range("K4").Select()
